# ImportUtil template update: add a "Form Tag" column (V) to the investor
# KYC import sheet, with sample values "Default" / "Gift City" for the
# first two sample investor rows, matching the new form_type tagging
# behaviour described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Form Tag" column
$ws.Range("V1").Value = "Form Tag"

# Sample values for the first two investor rows
$ws.Range("V2").Value = "Default"
$ws.Range("V3").Value = "Gift City"

# Move the view/selection over to the newly added column, as seen in the
# updated workbook (topLeftCell shifted to D1, active cell now V2)
$ws.Range("D1").Select()
$ws.Range("V2").Select()
